# Append a new data row (row 16) to the Swissmedic "galenic" package sheet,
# mirroring the structure of the existing rows (e.g. row 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

# Match the row height used by the other data rows (e.g. row 15).
$ws.Rows.Item($row).RowHeight = 12.75

$ws.Cells.Item($row, 1).Value  = 10386                         # Zulassungs-Nummer
$ws.Cells.Item($row, 2).Value  = 1                              # Dosisstärke-nummer
$ws.Cells.Item($row, 3).Value  = "Lapidar 4, Filmtabletten"     # Präparatebezeichnung
$ws.Cells.Item($row, 4).Value  = "Kräuterpfarrer Künzle AG"     # Zulassungsinhaberin
$ws.Cells.Item($row, 5).Value  = "02.08.1."                     # IT-Nummer
$ws.Cells.Item($row, 6).Value  = "C05CA51"                      # ATC-Code
$ws.Cells.Item($row, 7).Value  = "Synthetika human"             # Heilmittelcode
# Date columns: carry the same date number format used by the other rows
# (e.g. row 15) so the serials render as dates rather than plain numbers.
$dateFmt = $ws.Cells.Item(15, 8).NumberFormat
$ws.Cells.Item($row, 8).NumberFormat  = $dateFmt
$ws.Cells.Item($row, 8).Value  = 13027                          # Erstzul.datum Präp.
$ws.Cells.Item($row, 9).NumberFormat  = $dateFmt
$ws.Cells.Item($row, 9).Value  = 13027                          # Zul.datum Dosisstärke
$ws.Cells.Item($row, 10).NumberFormat = $dateFmt
$ws.Cells.Item($row, 10).Value = 42358                          # Gültigkeits-datum
$ws.Cells.Item($row, 11).Value = 47                             # Verpackungs ID

# Packungsgrösse is a textual "150" (not the number 150) in the source data,
# matching the other text columns even though it looks numeric -- force text
# via the number format so Excel doesn't auto-coerce it to a Number cell.
$ws.Cells.Item($row, 12).NumberFormat = "@"
$ws.Cells.Item($row, 12).Value = "150"
$ws.Cells.Item($row, 12).NumberFormat = "GENERAL"

$ws.Cells.Item($row, 13).Value = "Tablette(n)"                                                                  # Einheit
$ws.Cells.Item($row, 14).Value = "D"                                                                            # Abgabekategorie
$ws.Cells.Item($row, 15).Value = "rutosidum trihydricum, aescinum"                                              # Wirkstoff
$ws.Cells.Item($row, 16).Value = "rutosidum trihydricum 20 mg, aescinum 25 mg, aromatica, excipiens pro compresso." # Zusammensetzung
$ws.Cells.Item($row, 17).Value = "Symptome bei Krampfadern"                                                      # Anwendungsgebiet Präparate
